$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price and 1h volume change columns).
# D-column values are forced to text with a leading apostrophe so Excel
# does not reinterpret dotted price strings (e.g. "4.30") as numbers,
# matching the original inline-string cell type.

$ws.Range("D2").Value = "'26.648.55"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "'1.633.71"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'213.17"
$ws.Range("D6").Value = "'0.494"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").Value = "'0.0625"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  +2.71%  "
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").Value = "'1.860.97"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "'1.660.43"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").Value = "'26.651.43"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "'209.84"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'4.30"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "'147.14"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").Value = "'6.90"
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "'0.0522"
$ws.Range("E30").Value = "  +5.04%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").Value = "'2.96"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("E34").Value = "  +0.98%  "
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "'0.0172"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").Value = "'1.170.10"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Value = "'0.811"
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'0.506"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'0.794"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'1.770.49"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "'92.51"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'1.55"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").Value = "'54.70"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  -0.10%  "
